$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 'white athletic leggings mens'
    2 = 'supreme basketball tights'
    3 = 'snowmobiling knee pads'
    4 = 'photographer knee pads'
    5 = 'raymens leggings'
    6 = 'training basketball youth'
    7 = 'baleaf men compression pants'
    8 = 'winter leggings men'
    9 = 'kids pants with knee pads'
    10 = 'eclipse knee pads'
    11 = 'pilates knee pads'
    12 = 'fitted mens tights'
    13 = 'mens knee pads bmx'
    14 = 'canoe knee pad'
    15 = 'graduated compression pants'
    16 = 'graduated compression leggings'
    17 = 'soft knee pads'
    18 = 'knee pads nike'
    19 = 'knee pad toddler'
    20 = 'knee pad yoga'
    21 = 'knee pads dance'
    22 = 'knee pads mizuno'
    23 = 'elbow knee pads'
    24 = 'nee pads basketball'
    25 = 'compression pants set men'
    26 = 'mens basketball snap pants'
    27 = 'nike compression pants for youth'
    28 = 'mens compression tights cold weather'
    29 = 'mens under armour basketball tights'
    30 = 'mens tights pockets'
    31 = 'rollerblades knee pads'
    32 = 'knee pads 3xl'
    33 = 'knee pads 8'
    34 = 'knee pad hard'
    35 = 'men leggings fleece'
    36 = 'teflex knee pads'
    37 = 'man leggings thermal'
    38 = 'sailing knee pad'
    39 = 'knee pads downhill'
    40 = 'knee pads airsoft'
    41 = 'knee pads army'
    42 = 'knee pads enduro'
    43 = 'knee pads bike'
    44 = 'knee pads caving'
    45 = 'knee pads for women'
    46 = 'knee pads gloves'
    47 = 'knee pads green'
    48 = 'knee pads kali'
    49 = 'knee pads kuangmi'
    50 = 'knee pads longboard'
    51 = 'knee pads neoprene'
    52 = 'knee pads ocp'
    53 = 'knee pads orange'
    54 = 'knee pads over pants'
    55 = 'knee pads plastic'
    56 = 'knee pads purple'
    57 = 'knee pads razor'
    58 = 'knee pads red'
    59 = 'knee pads rollerblading'
    60 = 'knee pads sailing'
    61 = 'knee pads scooter'
    62 = 'knee pads set'
    63 = 'knee pads shooting'
    64 = 'knee pads swat'
    65 = 'knee pads teen'
    66 = 'knee pads tsg'
    67 = 'knee pads usmc'
    68 = 'knee pads viper'
    69 = 'knee pads white'
    70 = 'knee pads yellow'
    71 = 'mens leggings xs'
    72 = 'xtextile compression pants men'
    73 = 'crx men''s tights'
    74 = 'yoga capri pants'
    75 = 'dodoing kneepads'
    76 = 'cavaliers basketball leggings'
    77 = 'basketball knee pads kids'
    78 = 'knee pads for toddlers'
    79 = 'knee pads skating'
    80 = 'knee pads skateboarding'
    81 = 'basketball knee pads kids boys'
    82 = 'knee pads for dance'
    83 = 'knee pads rollerblade'
    84 = 'knee pads tan'
    85 = 'knee pad dancer'
    86 = 'knee pads adidas'
    87 = 'knee pads basketball mcdavid'
    88 = 'knee pads dancing'
    89 = 'knee pads dodgeball'
    90 = 'knee pads pair'
    91 = 'knee pads longboarding'
    92 = 'knee pads nba'
    93 = 'knee pads pole'
    94 = 'knee pad and elbow pads'
    95 = 'knee pad adidas'
    96 = 'knee pad asics'
    97 = 'knee pad for kids'
    98 = 'knee pad military'
    99 = 'knee pad mma'
    100 = 'knee pad mizuno'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
